# Add new match/round data rows (69-91) to the Jogadores sheet,
# following the same layout as the existing rows, and update the
# active window's scroll position / selection to reflect the new
# bottom of the data range (mirrors the manual Excel edit session).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Jogadores")

$ws.Cells.Item(69, 1).Value = 'Cabeleira'
$ws.Cells.Item(69, 3).Value = 1
$ws.Cells.Item(69, 4).Value = 2
$ws.Cells.Item(69, 5).Value = 3
$ws.Cells.Item(69, 6).Value = 1
$ws.Cells.Item(69, 7).Value = 1
$ws.Cells.Item(69, 8).Value = 0
$ws.Cells.Item(69, 9).Value = 1
$ws.Cells.Item(69, 10).Value = 0
$ws.Cells.Item(69, 11).Value = 0

$ws.Cells.Item(70, 1).Value = 'Peixe'
$ws.Cells.Item(70, 3).Value = 1
$ws.Cells.Item(70, 4).Value = 2
$ws.Cells.Item(70, 5).Value = 3
$ws.Cells.Item(70, 6).Value = 0
$ws.Cells.Item(70, 7).Value = 1
$ws.Cells.Item(70, 8).Value = 0
$ws.Cells.Item(70, 9).Value = 1
$ws.Cells.Item(70, 10).Value = 0
$ws.Cells.Item(70, 11).Value = 0

$ws.Cells.Item(71, 1).Value = 'Leandro'
$ws.Cells.Item(71, 3).Value = 1
$ws.Cells.Item(71, 4).Value = 2
$ws.Cells.Item(71, 5).Value = 3
$ws.Cells.Item(71, 6).Value = 0
$ws.Cells.Item(71, 7).Value = 1
$ws.Cells.Item(71, 8).Value = 0
$ws.Cells.Item(71, 9).Value = 1
$ws.Cells.Item(71, 10).Value = 0
$ws.Cells.Item(71, 11).Value = 0

$ws.Cells.Item(72, 1).Value = 'Jorge'
$ws.Cells.Item(72, 3).Value = 1
$ws.Cells.Item(72, 4).Value = 2
$ws.Cells.Item(72, 5).Value = 3
$ws.Cells.Item(72, 6).Value = 0
$ws.Cells.Item(72, 7).Value = 1
$ws.Cells.Item(72, 8).Value = 0
$ws.Cells.Item(72, 9).Value = 1
$ws.Cells.Item(72, 10).Value = 0
$ws.Cells.Item(72, 11).Value = 0

$ws.Cells.Item(73, 1).Value = 'Philipe'
$ws.Cells.Item(73, 3).Value = 1
$ws.Cells.Item(73, 4).Value = 2
$ws.Cells.Item(73, 5).Value = 3
$ws.Cells.Item(73, 6).Value = 2
$ws.Cells.Item(73, 7).Value = 1
$ws.Cells.Item(73, 8).Value = 0
$ws.Cells.Item(73, 9).Value = 1
$ws.Cells.Item(73, 10).Value = 0
$ws.Cells.Item(73, 11).Value = 0

$ws.Cells.Item(74, 1).Value = 'Joazinho'
$ws.Cells.Item(74, 3).Value = 2
$ws.Cells.Item(74, 4).Value = 1
$ws.Cells.Item(74, 5).Value = 3
$ws.Cells.Item(74, 6).Value = 1
$ws.Cells.Item(74, 7).Value = 1
$ws.Cells.Item(74, 8).Value = 0
$ws.Cells.Item(74, 9).Value = 0
$ws.Cells.Item(74, 10).Value = 0
$ws.Cells.Item(74, 11).Value = 0

$ws.Cells.Item(75, 1).Value = 'Marcelão'
$ws.Cells.Item(75, 3).Value = 2
$ws.Cells.Item(75, 4).Value = 1
$ws.Cells.Item(75, 5).Value = 3
$ws.Cells.Item(75, 6).Value = 2
$ws.Cells.Item(75, 7).Value = 1
$ws.Cells.Item(75, 8).Value = 0
$ws.Cells.Item(75, 9).Value = 0
$ws.Cells.Item(75, 10).Value = 0
$ws.Cells.Item(75, 11).Value = 0

$ws.Cells.Item(76, 1).Value = 'Leandrao'
$ws.Cells.Item(76, 3).Value = 2
$ws.Cells.Item(76, 4).Value = 1
$ws.Cells.Item(76, 5).Value = 3
$ws.Cells.Item(76, 6).Value = 1
$ws.Cells.Item(76, 7).Value = 1
$ws.Cells.Item(76, 8).Value = 0
$ws.Cells.Item(76, 9).Value = 0
$ws.Cells.Item(76, 10).Value = 0
$ws.Cells.Item(76, 11).Value = 0

$ws.Cells.Item(77, 1).Value = 'David'
$ws.Cells.Item(77, 3).Value = 2
$ws.Cells.Item(77, 4).Value = 1
$ws.Cells.Item(77, 5).Value = 3
$ws.Cells.Item(77, 6).Value = 0
$ws.Cells.Item(77, 7).Value = 1
$ws.Cells.Item(77, 8).Value = 0
$ws.Cells.Item(77, 9).Value = 0
$ws.Cells.Item(77, 10).Value = 0
$ws.Cells.Item(77, 11).Value = 0

$ws.Cells.Item(78, 1).Value = 'Marcos'
$ws.Cells.Item(78, 3).Value = 2
$ws.Cells.Item(78, 4).Value = 1
$ws.Cells.Item(78, 5).Value = 3
$ws.Cells.Item(78, 6).Value = 2
$ws.Cells.Item(78, 7).Value = 1
$ws.Cells.Item(78, 8).Value = 0
$ws.Cells.Item(78, 9).Value = 0
$ws.Cells.Item(78, 10).Value = 0
$ws.Cells.Item(78, 11).Value = 0

$ws.Cells.Item(79, 1).Value = 'Juscielio'
$ws.Cells.Item(79, 3).Value = 4
$ws.Cells.Item(79, 4).Value = 3
$ws.Cells.Item(79, 5).Value = 0
$ws.Cells.Item(79, 6).Value = 0
$ws.Cells.Item(79, 7).Value = 1
$ws.Cells.Item(79, 8).Value = 1
$ws.Cells.Item(79, 9).Value = 0
$ws.Cells.Item(79, 10).Value = 1
$ws.Cells.Item(79, 11).Value = 0

$ws.Cells.Item(80, 1).Value = 'Guinha'
$ws.Cells.Item(80, 3).Value = 4
$ws.Cells.Item(80, 4).Value = 3
$ws.Cells.Item(80, 5).Value = 0
$ws.Cells.Item(80, 6).Value = 1
$ws.Cells.Item(80, 7).Value = 1
$ws.Cells.Item(80, 8).Value = 1
$ws.Cells.Item(80, 9).Value = 0
$ws.Cells.Item(80, 10).Value = 0
$ws.Cells.Item(80, 11).Value = 0

$ws.Cells.Item(81, 1).Value = 'Euzebio'
$ws.Cells.Item(81, 3).Value = 4
$ws.Cells.Item(81, 4).Value = 3
$ws.Cells.Item(81, 5).Value = 0
$ws.Cells.Item(81, 6).Value = 1
$ws.Cells.Item(81, 7).Value = 1
$ws.Cells.Item(81, 8).Value = 1
$ws.Cells.Item(81, 9).Value = 0
$ws.Cells.Item(81, 10).Value = 0
$ws.Cells.Item(81, 11).Value = 0

$ws.Cells.Item(82, 1).Value = 'Eder'
$ws.Cells.Item(82, 3).Value = 4
$ws.Cells.Item(82, 4).Value = 3
$ws.Cells.Item(82, 5).Value = 0
$ws.Cells.Item(82, 6).Value = 1
$ws.Cells.Item(82, 7).Value = 1
$ws.Cells.Item(82, 8).Value = 1
$ws.Cells.Item(82, 9).Value = 0
$ws.Cells.Item(82, 10).Value = 0
$ws.Cells.Item(82, 11).Value = 0

$ws.Cells.Item(83, 1).Value = 'Leandrinho'
$ws.Cells.Item(83, 3).Value = 4
$ws.Cells.Item(83, 4).Value = 3
$ws.Cells.Item(83, 5).Value = 0
$ws.Cells.Item(83, 6).Value = 2
$ws.Cells.Item(83, 7).Value = 1
$ws.Cells.Item(83, 8).Value = 1
$ws.Cells.Item(83, 9).Value = 0
$ws.Cells.Item(83, 10).Value = 0
$ws.Cells.Item(83, 11).Value = 0

$ws.Cells.Item(84, 1).Value = 'Ismael'
$ws.Cells.Item(84, 3).Value = 1
$ws.Cells.Item(84, 4).Value = 2
$ws.Cells.Item(84, 5).Value = 2
$ws.Cells.Item(84, 6).Value = 1
$ws.Cells.Item(84, 7).Value = 1
$ws.Cells.Item(84, 8).Value = 0
$ws.Cells.Item(84, 9).Value = 0
$ws.Cells.Item(84, 10).Value = 0
$ws.Cells.Item(84, 11).Value = 0

$ws.Cells.Item(85, 1).Value = 'Boneco'
$ws.Cells.Item(85, 3).Value = 1
$ws.Cells.Item(85, 4).Value = 2
$ws.Cells.Item(85, 5).Value = 2
$ws.Cells.Item(85, 6).Value = 0
$ws.Cells.Item(85, 7).Value = 1
$ws.Cells.Item(85, 8).Value = 0
$ws.Cells.Item(85, 9).Value = 0
$ws.Cells.Item(85, 10).Value = 0
$ws.Cells.Item(85, 11).Value = 0

$ws.Cells.Item(86, 1).Value = 'Corinthiano'
$ws.Cells.Item(86, 3).Value = 1
$ws.Cells.Item(86, 4).Value = 2
$ws.Cells.Item(86, 5).Value = 2
$ws.Cells.Item(86, 6).Value = 0
$ws.Cells.Item(86, 7).Value = 1
$ws.Cells.Item(86, 8).Value = 0
$ws.Cells.Item(86, 9).Value = 0
$ws.Cells.Item(86, 10).Value = 0
$ws.Cells.Item(86, 11).Value = 0

$ws.Cells.Item(87, 1).Value = 'Athos'
$ws.Cells.Item(87, 3).Value = 1
$ws.Cells.Item(87, 4).Value = 2
$ws.Cells.Item(87, 5).Value = 2
$ws.Cells.Item(87, 6).Value = 0
$ws.Cells.Item(87, 7).Value = 1
$ws.Cells.Item(87, 8).Value = 0
$ws.Cells.Item(87, 9).Value = 0
$ws.Cells.Item(87, 10).Value = 0
$ws.Cells.Item(87, 11).Value = 0

$ws.Cells.Item(88, 1).Value = 'Ranyeri'
$ws.Cells.Item(88, 3).Value = 1
$ws.Cells.Item(88, 4).Value = 2
$ws.Cells.Item(88, 5).Value = 2
$ws.Cells.Item(88, 6).Value = 1
$ws.Cells.Item(88, 7).Value = 1
$ws.Cells.Item(88, 8).Value = 0
$ws.Cells.Item(88, 9).Value = 0
$ws.Cells.Item(88, 10).Value = 0
$ws.Cells.Item(88, 11).Value = 0

$ws.Cells.Item(89, 1).Value = 'Matheus'
$ws.Cells.Item(89, 3).Value = 4
$ws.Cells.Item(89, 4).Value = 4
$ws.Cells.Item(89, 5).Value = 1
$ws.Cells.Item(89, 6).Value = 0
$ws.Cells.Item(89, 7).Value = 1
$ws.Cells.Item(89, 8).Value = 0
$ws.Cells.Item(89, 9).Value = 0
$ws.Cells.Item(89, 10).Value = 0
$ws.Cells.Item(89, 11).Value = 3

$ws.Cells.Item(90, 1).Value = 'Lucian'
$ws.Cells.Item(90, 3).Value = 1
$ws.Cells.Item(90, 4).Value = 1
$ws.Cells.Item(90, 5).Value = 5
$ws.Cells.Item(90, 6).Value = 0
$ws.Cells.Item(90, 7).Value = 1
$ws.Cells.Item(90, 8).Value = 0
$ws.Cells.Item(90, 9).Value = 0
$ws.Cells.Item(90, 10).Value = 0
$ws.Cells.Item(90, 11).Value = 10

$ws.Cells.Item(91, 1).Value = 'Chelin'
$ws.Cells.Item(91, 3).Value = 3
$ws.Cells.Item(91, 4).Value = 3
$ws.Cells.Item(91, 5).Value = 2
$ws.Cells.Item(91, 6).Value = 0
$ws.Cells.Item(91, 7).Value = 1
$ws.Cells.Item(91, 8).Value = 0
$ws.Cells.Item(91, 9).Value = 0
$ws.Cells.Item(91, 10).Value = 0
$ws.Cells.Item(91, 11).Value = 4

# Move the frozen pane's scroll position down to the new rows and
# select the first empty cell below the newly entered data, matching
# where the user's cursor ended up after typing the last row.
$win = $excel.ActiveWindow
$win.ScrollRow = 69
$win.ScrollColumn = 1
$ws.Range("A92").Select()
